# This script rotates the contents of columns D, E and F (category-code,
# group-name, category-name) for every row of the SectorGroup sheet:
#   new D (category-code) = old F
#   new E (group-name)    = old D
#   new F (category-name) = old E
# Column G (group-code) and columns A-C are left untouched.
# This mirrors the column header re-ordering:
#   D,E,F,G : group-name, category-name, category-code, group-code
#        -> category-code, group-name, category-name, group-code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Read columns D, E, F into arrays first so the row-by-row write doesn't
# clobber values that still need to be read.
$dVals = @()
$eVals = @()
$fVals = @()

for ($r = 1; $r -le $lastRow; $r++) {
    $dVals += , ($ws.Cells.Item($r, 4).Value2)
    $eVals += , ($ws.Cells.Item($r, 5).Value2)
    $fVals += , ($ws.Cells.Item($r, 6).Value2)
}

for ($r = 1; $r -le $lastRow; $r++) {
    $idx = $r - 1
    $ws.Cells.Item($r, 4).Value2 = $fVals[$idx]
    $ws.Cells.Item($r, 5).Value2 = $dVals[$idx]
    $ws.Cells.Item($r, 6).Value2 = $eVals[$idx]
}
